# Generate Report for Handback
# Applies the handback-report update to the Overview / zh-cn / de-de sheets:
#  - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#    (every cell that shows the status needs to be touched, since each
#    worksheet cell owns its own shared-string slot)
#  - zh-cn / de-de rows gain a "Latest Target File" hyperlink (col I) and a
#    "Latest Handback File" name (col J), plus a refreshed
#    "Latest Handback DateTime" (col K)
#  - A handful of columns get widened to fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$baseFileName = "2cb00a8e-47d0-4f03-9380-58ea0d335076.md"
$githubUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1582352ed706b9c68b78a32364222b694b12733e/e2e/2cb00a8e-47d0-4f03-9380-58ea0d335076.md"
$statusText   = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Status column -> "Handed back: in sync with en-US" everywhere it shows
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsZhCn.Range("C2").Value = $statusText
$wsDeDe.Range("C2").Value = $statusText

# ---------------------------------------------------------------------------
# 2. zh-cn row: Latest Target File (I2) / Latest Handback File (J2) /
#    Latest Handback DateTime (K2)
# ---------------------------------------------------------------------------
$wsZhCn.Range("I2").Value = $baseFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $githubUrl, "", "", $baseFileName)
$wsZhCn.Range("J2").Value = "2cb00a8e-47d0-4f03-9380-58ea0d335076.04182e2239ef5e8d874a9375d0630731b31df6e1.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-09 13:04:15"

# ---------------------------------------------------------------------------
# 3. de-de row: Latest Target File (I2) / Latest Handback File (J2) /
#    Latest Handback DateTime (K2)
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value = $baseFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $githubUrl, "", "", $baseFileName)
$wsDeDe.Range("J2").Value = "2cb00a8e-47d0-4f03-9380-58ea0d335076.04182e2239ef5e8d874a9375d0630731b31df6e1.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-09 13:04:33"

# ---------------------------------------------------------------------------
# 4. Column widths - widen to fit the longer handback content
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.15   # zh-cn status column
$wsOverview.Columns.Item(6).ColumnWidth = 29.15   # de-de status column

$wsZhCn.Columns.Item(3).ColumnWidth = 29.15    # Status
$wsZhCn.Columns.Item(9).ColumnWidth = 39.15    # Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.15   # Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = 29.15    # Status
$wsDeDe.Columns.Item(9).ColumnWidth = 39.15    # Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.15   # Latest Handback File
